$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" (D) and "Volume(1h)" (E) columns store plain text values
# (the source site renders numbers like "43.825.39" or "  -0.28%  " as
# text). Force each touched cell to the Text number format first so
# Excel does not reinterpret a numeric-looking string (e.g. "72.59")
# as a real number and silently change its stored representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.825.39"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.343.45"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "239.03"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("D7").Value = "72.59"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -7.38%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "58.79"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "32.73"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "7.24"
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("D15").Value = "2.693.66"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "16.07"
$ws.Range("E16").Value = "  -4.50%  "
$ws.Range("D18").Value = "2.338.84"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "43.763.20"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "78.25"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "251.24"
$ws.Range("E23").Value = "  -4.36%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("E28").Value = "  -5.43%  "
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "176.59"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "22.20"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("D36").Value = "5.34"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "3.72"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "5.58"
$ws.Range("E40").Value = "  +21.77%  "
$ws.Range("D41").Value = "0.0270"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").Value = "65.27"
$ws.Range("E42").Value = "  +16.06%  "
$ws.Range("D43").Value = "9.22"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").Value = "0.105"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "18.76"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("E46").Value = "  -10.34%  "
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.40"
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "97.91"
$ws.Range("E51").Value = "  -4.09%  "
